# Lesson 1.3 The Data Design Recipe.pptx
# 1) Slide 21 ("The DDR for itemization data"): reword the trailing
#    sentence of the highlighted TextBox so it reads "...select the
#    relevant case." instead of "...select each case."
# 2) Slide 7 ("Output of DDR Step 1"): swap "mixed" for "itemization"
#    in the list of data kinds.

$p = $ppt.ActivePresentation

# --- Slide 21: TextBox 5 -------------------------------------------------
$slide21 = $p.Slides.Item(21)
$box = $slide21.Shapes.Item(4)
$tr = $box.TextFrame.TextRange

$fullText = $tr.Text
$oldPhrase = "each case"
$newPhrase = "the relevant case"
$pos = $fullText.IndexOf($oldPhrase)
if ($pos -ge 0) {
    $target = $tr.Characters($pos + 1, $oldPhrase.Length)
    $target.Text = $newPhrase
}

# --- Slide 7: Content Placeholder 2, paragraph 1 -------------------------
$slide7 = $p.Slides.Item(7)
$content = $slide7.Shapes.Item(2)
$ctr = $content.TextFrame.TextRange

$para1 = $ctr.Paragraphs(1, 1)
$para1Text = $para1.Text.TrimEnd([char]13)
$updated = $para1Text.Replace("mixed", "itemization")
$paraRange = $ctr.Characters($para1.Start, $para1Text.Length)
$paraRange.Text = $updated
